$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57-108 down to 58-109
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new record
$ws.Range("A57").Value2 = 9
$ws.Range("B57").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C57").Value2 = "Metropolitana"
$ws.Range("D57").Value2 = 44825
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E57").Value2 = 13
$ws.Range("F57").Value2 = 100112005
$ws.Range("G57").Value2 = "Puerro"
$ws.Range("H57").Value2 = "Sin especificar"
$ws.Range("I57").Value2 = "Primera"
$ws.Range("J57").Value2 = 70
$ws.Range("K57").Value2 = 12000
$ws.Range("L57").Value2 = 12000
$ws.Range("M57").Value2 = 12000
$ws.Range("N57").Value2 = "$/paquete 20 unidades"
$ws.Range("O57").Value2 = "Provincia de Chacabuco"
$ws.Range("P57").Value2 = 600
$ws.Range("Q57").Value2 = 20
$ws.Range("R57").Value2 = "Hortaliza"

# Append a brand new row 110 at the end of the sheet
$ws.Range("A110").Value2 = 9
$ws.Range("B110").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C110").Value2 = "Metropolitana"
$ws.Range("D110").Value2 = 44832
$ws.Range("D110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E110").Value2 = 13
$ws.Range("F110").Value2 = 100112005
$ws.Range("G110").Value2 = "Puerro"
$ws.Range("H110").Value2 = "Sin especificar"
$ws.Range("I110").Value2 = "Segunda"
$ws.Range("J110").Value2 = 30
$ws.Range("K110").Value2 = 12000
$ws.Range("L110").Value2 = 12000
$ws.Range("M110").Value2 = 12000
$ws.Range("N110").Value2 = "$/paquete 20 unidades"
$ws.Range("O110").Value2 = "Provincia de Melipilla"
$ws.Range("P110").Value2 = 600
$ws.Range("Q110").Value2 = 20
$ws.Range("R110").Value2 = "Hortaliza"
